$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45734 -> 45735, i.e. 2025-03-18 -> 2025-03-19) for every data row (2..44).
for ($row = 2; $row -le 44; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45734) {
        $cell.Value2 = 45735
    }
}
